{"js": "// Update the Table 3 survival-probability figures for age classes\n// G (4-6y) and H (6-8y), and two cells in I (>8y), per the authored diff.\n//\n// The table (first/only table in the document) has columns:\n//   0: Age Class, 1: Benkovac-Barice MN, 2: Islam Gr\u010dki MN,\n//   3: Smil\u010di\u0107 EN, 4: Smil\u010di\u0107 MN, 5: Zemunik Donji MN\n// and rows (0-indexed, row 0 = header):\n//   7: G (4-6y), 8: H (6-8y), 9: I (>8y)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body.\");\n}\n\nconst table = tables.items[0];\n\n// List of edits: [rowIndex, colIndex, oldText, newText]\nconst edits = [\n  [7, 1, \"3.24\", \"1.63\"],\n  [7, 2, \"9.38\", \"4.80\"],\n  [7, 3, \"2.27\", \"1.14\"],\n  [7, 4, \"9.65\", \"4.95\"],\n  [8, 1, \"0.00\", \"1.63\"],\n  [8, 2, \"6.25\", \"4.80\"],\n  [8, 3, \"0.00\", \"1.14\"],\n  [8, 4, \"2.19\", \"4.95\"],\n  [9, 2, \"0.00\", \"3.18\"],\n  [9, 4, \"0.00\", \"1.10\"],\n];\n\n// Load every target cell's current value first so we can validate before\n// writing (defensive: fail loudly instead of silently mis-editing).\nconst cells = edits.map(([r, c]) => table.getCell(r, c));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nedits.forEach(([r, c, oldText, newText], i) => {\n  const cell = cells[i];\n  const current = (cell.value || \"\").trim();\n  if (current !== oldText) {\n    throw new Error(\n      `Unexpected value at row ${r}, col ${c}: expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n  cell.value = newText;\n});\n\nawait context.sync();\n", "ps1": "# Update the Table 3 survival-probability figures for age classes\n# G (4-6y) and H (6-8y), and two cells in I (>8y), per the authored diff.\n#\n# Table 1 (the only table in the document) columns are:\n#   1: Age Class, 2: Benkovac-Barice MN, 3: Islam Gr\u010dki MN,\n#   4: Smil\u010di\u0107 EN, 5: Smil\u010di\u0107 MN, 6: Zemunik Donji MN\n# Rows (1-indexed, row 1 = header):\n#   8: G (4-6y), 9: H (6-8y), 10: I (>8y)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Each entry: row, column, expected current text, new text\n$edits = @(\n  @(8, 2, \"3.24\", \"1.63\"),\n  @(8, 3, \"9.38\", \"4.80\"),\n  @(8, 4, \"2.27\", \"1.14\"),\n  @(8, 5, \"9.65\", \"4.95\"),\n  @(9, 2, \"0.00\", \"1.63\"),\n  @(9, 3, \"6.25\", \"4.80\"),\n  @(9, 4, \"0.00\", \"1.14\"),\n  @(9, 5, \"2.19\", \"4.95\"),\n  @(10, 3, \"0.00\", \"3.18\"),\n  @(10, 5, \"0.00\", \"1.10\")\n)\n\nforeach ($edit in $edits) {\n  $row = $edit[0]\n  $col = $edit[1]\n  $oldText = $edit[2]\n  $newText = $edit[3]\n\n  $cell = $table.Cell($row, $col)\n  $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n  if ($current -ne $oldText) {\n    throw \"Unexpected value at row ${row}, col ${col}: expected '$oldText' but found '$current'\"\n  }\n\n  $cell.Range.Text = $newText\n}\n"}
